$d = $word.ActiveDocument

# 1. Update the Java download URL (step 1 instructions).
$r1 = $d.Content
$found1 = $r1.Find.Execute(
    "1. Go to https://www.oracle.com/java/technologies/downloads/#java8-windows",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1. Go to https://www.oracle.com/uk/java/technologies/javase/javase8-archive-downloads.html",
    2)
Write-Output "URL replace found: $found1"

# 2. Fix the "7. Download settings.xml from " text, which was previously
#    split across three runs ("7. Download se" + "t" + "tings.xml from ").
#    Replacing the full phrase coalesces it into a single run/string.
$r2 = $d.Content
$found2 = $r2.Find.Execute(
    "7. Download settings.xml from ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "7. Download settings.xml from ",
    2)
Write-Output "Settings text replace found: $found2"
